# "Addition of Demo xml" - refresh the demo data shown on the
# "Sec invoice Master" sheet with a new sample record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify the header for the "Date CREATED" column by spelling out the
# expected date format.
$ws.Range("B1").Value = "Date CREATED (MM/DD/YYYY)"

# Refresh the sample/demo row with a new FC Order ID, tracking number and
# secondary invoice number.
$ws.Range("A2").Value = "57466645"
$ws.Range("C2").Value = "FCT877055778467348480"
$ws.Range("I2").Value = "57466645+1"

# Reflect the user's latest view/selection state on the sheet.
$ws.Range("L1").Select() | Out-Null
